# Fixed typo in image
#
# 1) The "Date Placeholder" field on the slide master and every slide
#    layout shows the auto date/time field for the deck; bump it from
#    12/2/2016 to 12/7/2016 everywhere it appears.
# 2) Fix a copy/paste typo on the pinout diagram: the pin labelled
#    "P10" should read "P9".

$p = $ppt.ActivePresentation

function Update-DatePlaceholders($shapes) {
    for ($idx = 1; $idx -le $shapes.Count; $idx++) {
        $shp = $shapes.Item($idx)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "12/7/2016"
        }
    }
}

# Slide master's own date placeholder.
Update-DatePlaceholders $p.SlideMaster.Shapes

# Every layout hanging off the master has its own date placeholder too.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes
}

# Fix the "P10" -> "P9" typo on the pinout image (slide 1, inside the
# top-level "Group 1" group of shapes).
$slide = $p.Slides.Item(1)
$group = $slide.Shapes.Item(1)
for ($gi = 1; $gi -le $group.GroupItems.Count; $gi++) {
    $shp = $group.GroupItems.Item($gi)
    if ($shp.Name -eq "Rounded Rectangle 96") {
        $shp.TextFrame.TextRange.Text = "P9"
    }
}
